$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.479846835136414
$ws.Range("B1").Value = 2.022845506668091
$ws.Range("C1").Value = 3.124520301818848
$ws.Range("D1").Value = 4.745704650878906
$ws.Range("E1").Value = 0.8732813000679016
